# Update cryptos list: refreshed prices / 1h volume %, and
# RenderToken / Stacks swapped position (rows 37-38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "67.893.12"
    "E2" = "  -1.28%  "
    "D3" = "2.401.10"
    "E3" = "  -2.21%  "
    "D4" = "1.00"
    "E4" = "  -0.01%  "
    "D5" = "552.01"
    "E5" = "  -1.12%  "
    "D6" = "157.72"
    "E6" = "  -3.02%  "
    "E7" = "  +0.00%  "
    "D8" = "0.503"
    "E8" = "  -0.01%  "
    "E9" = "  +4.12%  "
    "E10" = "  -1.42%  "
    "D11" = "0.326"
    "E11" = "  -2.43%  "
    "D12" = "4.70"
    "E12" = "  -2.38%  "
    "D13" = "67.801.43"
    "E13" = "  -1.25%  "
    "E14" = "  -0.28%  "
    "D15" = "22.76"
    "E15" = "  -3.38%  "
    "D16" = "10.25"
    "E16" = "  -4.82%  "
    "D17" = "328.68"
    "E17" = "  -3.63%  "
    "D18" = "6.77"
    "E18" = "  -4.21%  "
    "D19" = "3.75"
    "E19" = "  -1.04%  "
    "E20" = "  -0.04%  "
    "D21" = "1.86"
    "E21" = "  -4.05%  "
    "D22" = "65.55"
    "E22" = "  -2.07%  "
    "D23" = "3.60"
    "E23" = "  -2.48%  "
    "D24" = "8.01"
    "E24" = "  -2.25%  "
    "D25" = "0.0₃0791"
    "E25" = "  -3.22%  "
    "D26" = "7.02"
    "E26" = "  -1.90%  "
    "D27" = "1.00"
    "E27" = "  +0.03%  "
    "D28" = "418.06"
    "D29" = "1.13"
    "E29" = "  -1.72%  "
    "E30" = "  -1.99%  "
    "D31" = "157.10"
    "E31" = "  -0.23%  "
    "D32" = "18.98"
    "E32" = "  -0.20%  "
    "E33" = "  -0.10%  "
    "D34" = "17.62"
    "E34" = "  -0.86%  "
    "E35" = "  -3.97%  "
    "E36" = "  -3.19%  "
    "B37" = "Stacks"
    "C37" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D37" = "1.45"
    "E37" = "  -1.54%  "
    "B38" = "RenderToken"
    "C38" = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
    "D38" = "4.21"
    "E38" = "  -5.55%  "
    "E39" = "  -3.99%  "
    "E40" = "  -2.04%  "
    "D41" = "128.54"
    "E41" = "  -3.34%  "
    "D42" = "1.92"
    "E42" = "  -7.62%  "
    "D43" = "0.0703"
    "E43" = "  -2.00%  "
    "D44" = "0.473"
    "E44" = "  -2.02%  "
    "D45" = "0.554"
    "E45" = "  -0.95%  "
    "D46" = "0.0911"
    "E46" = "  +0.51%  "
    "E47" = "  -0.65%  "
    "E48" = "  -6.79%  "
    "D49" = "16.32"
    "E49" = "  -3.40%  "
    "E50" = "  -8.08%  "
    "E51" = "  -0.61%  "
}

foreach ($ref in $updates.Keys) {
    # Force text storage so numeric-looking strings (e.g. "1.00", "4.70")
    # keep their exact formatting instead of being coerced to a Double.
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
